$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current header (row1) and data (row2) values for columns C:G
# before they get overwritten, using Value2 to read real scalars.
$headers = @{}
$values = @{}
foreach ($col in @("C","D","E","F","G")) {
    $headers[$col] = $ws.Range("$($col)1").Value2
    $values[$col]  = $ws.Range("$($col)2").Value2
}

# Reorder the site columns C:G into the new order: PKV, STL, THL, AZC, SDU
# Target column <- source column (mapped by the original AZC/PKV/SDU/STL/THL data)
$ws.Range("C1").Value = $headers["D"]   # PKV
$ws.Range("C2").Value = $values["D"]

$ws.Range("D1").Value = $headers["F"]   # STL
$ws.Range("D2").Value = $values["F"]

$ws.Range("E1").Value = $headers["G"]   # THL
$ws.Range("E2").Value = $values["G"]

$ws.Range("F1").Value = $headers["C"]   # AZC
$ws.Range("F2").Value = $values["C"]

$ws.Range("G1").Value = $headers["E"]   # SDU
$ws.Range("G2").Value = $values["E"]

# Add new "Date" header in A1 (plain/default styling, no border/bold)
$ws.Range("A1").Value = "Date"

# Update A2 date value: drop the leading underscore ("_20240226" -> "20240226"),
# keeping it stored as text (format the cell as Text first so Excel does not
# auto-convert the digit string into a number).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "20240226"
